# project_list.xlsx - switch the "Collection mode" column from "auto" to
# "manual" on the "Project list" sheet (column AT), and update the sheet's
# active selection to match where the author left the cursor (AT13:AT14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project list")

# Column AT ("Collection mode") holds "auto" for both data rows (2 and 3);
# flip them to "manual".
$ws.Range("AT2").Value = "manual"
$ws.Range("AT3").Value = "manual"

# Make sure this sheet is the active one, then move the selection to
# AT13:AT14 (active cell AT13), matching the saved view state.
$ws.Activate()
$ws.Range("AT13:AT14").Select()
